# Rerun and summarise models without urban landuse
# 1) Rename all 9 summary sheets to their new "summ<id>" names
# 2) Update the "Education[T.Unknown]" label to "Education[T.Unknown/Other]"
#    in cell A5 on every sheet

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ52806160",
    "summ53047536",
    "summ53333990",
    "summ53643605",
    "summ53960419",
    "summ54246088",
    "summ54517683",
    "summ54786780",
    "summ55096220"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]

    if ($ws.Range("A5").Value2 -eq "Education[T.Unknown]") {
        $ws.Range("A5").Value = "Education[T.Unknown/Other]"
    }
}
